$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new "missing_values" worksheet after the last existing sheet ("5d")
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "missing_values"

# Header row
$ws5.Range("A1").Value = "age"
$ws5.Range("B1").Value = "sex\time"
$ws5.Range("C1").Value = 2007
$ws5.Range("D1").Value = 2010
$ws5.Range("E1").Value = 2013

# Data rows (same table as sheet "3d" but with two rows removed to
# represent missing combinations: age=1/sex=H and age=4/sex=F)
$ws5.Range("A2").Value = 0
$ws5.Range("B2").Value = "F"
$ws5.Range("C2").Value = 3722
$ws5.Range("D2").Value = 3395
$ws5.Range("E2").Value = 3347

$ws5.Range("A3").Value = 0
$ws5.Range("B3").Value = "H"
$ws5.Range("C3").Value = 338
$ws5.Range("D3").Value = 316
$ws5.Range("E3").Value = 323

$ws5.Range("A4").Value = 1
$ws5.Range("B4").Value = "F"
$ws5.Range("C4").Value = 2878
$ws5.Range("D4").Value = 2791
$ws5.Range("E4").Value = 2822

$ws5.Range("A5").Value = 2
$ws5.Range("B5").Value = "F"
$ws5.Range("C5").Value = 4073
$ws5.Range("D5").Value = 4161
$ws5.Range("E5").Value = 4429

$ws5.Range("A6").Value = 2
$ws5.Range("B6").Value = "H"
$ws5.Range("C6").Value = 1561
$ws5.Range("D6").Value = 1463
$ws5.Range("E6").Value = 1467

$ws5.Range("A7").Value = 3
$ws5.Range("B7").Value = "F"
$ws5.Range("C7").Value = 3507
$ws5.Range("D7").Value = 3741
$ws5.Range("E7").Value = 3366

$ws5.Range("A8").Value = 3
$ws5.Range("B8").Value = "H"
$ws5.Range("C8").Value = 2052
$ws5.Range("D8").Value = 2052
$ws5.Range("E8").Value = 2118

$ws5.Range("A9").Value = 4
$ws5.Range("B9").Value = "H"
$ws5.Range("C9").Value = 3785
$ws5.Range("D9").Value = 3508
$ws5.Range("E9").Value = 3172

# The new sheet becomes the active / tab-selected sheet.
$ws5.Select() | Out-Null
$ws5.Range("G24").Select() | Out-Null

# ---------------------------------------------------------------------------
# Selections left behind on "2d" and "3d" from browsing the workbook.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2d")
$ws2.Select() | Out-Null
$ws2.Range("C8").Select() | Out-Null

$ws3 = $wb.Worksheets.Item("3d")
$ws3.Select() | Out-Null
$ws3.Range("B15").Select() | Out-Null

# Re-select the new sheet last so it ends up as the active tab.
$ws5.Select() | Out-Null
